$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.182.12"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.901.53"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.28"
$ws.Range("E5").Value = "  +3.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.698"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.06"
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("E9").Value = "  +3.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.89"
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0755"
$ws.Range("E11").Value = "  +4.11%  "
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.177.66"
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "13.01"
$ws.Range("E14").Value = "  +5.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.738"
$ws.Range("E15").Value = "  +3.56%  "
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("D17").Value = "1.911.59"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "35.166.92"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.78"
$ws.Range("E19").Value = "  +2.07%  "
$ws.Range("D20").Value = "0.0₃0835"
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "242.89"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.00"
$ws.Range("E22").Value = "  +3.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.06"
$ws.Range("E23").Value = "  +4.80%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  +4.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.29"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.82"
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.60"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.51"
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("D31").Value = "4.127.97"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("E32").Value = "  +16.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0609"
$ws.Range("E33").Value = "  +6.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.32"
$ws.Range("E34").Value = "  +3.67%  "
$ws.Range("E35").Value = "  +10.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.23"
$ws.Range("E36").Value = "  +2.41%  "
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("E38").Value = "  -11.43%  "
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "102.05"
$ws.Range("E40").Value = "  +12.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.19"
$ws.Range("E41").Value = "  +6.65%  "
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0647"
$ws.Range("E44").Value = "  -5.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.44"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").Value = "1.319.63"
$ws.Range("E46").Value = "  -1.84%  "
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("E48").Value = "  -1.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.62"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.79"
$ws.Range("E50").Value = "  -6.38%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.03"
$ws.Range("E51").Value = "  -7.90%  "
